# Switch to version v1.18 of the Excel file.
#
# 1. Make "5. Collections" (2nd sheet) the active tab instead of
#    "1. Content items" (1st sheet).
# 2. On "5. Collections", the data block that lived in columns S:AF
#    (with columns A:Q hidden/blank) moves to columns A:N - i.e. the
#    18 leading (hidden) columns are removed entirely.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Move the data on "5. Collections" from S:AF to A:N -------------------
# Columns A:R (1-18) in front of the real data are removed, shifting
# everything 18 columns to the left and dropping the old hidden/blank
# columns entirely.
$ws2.Range("A1:R9").EntireColumn.Delete()

# --- Make "5. Collections" the active/selected sheet ----------------------
# Activating it updates workbookView/activeTab, flips tabSelected from
# sheet1's sheetView to sheet2's sheetView, and clears the old
# topLeftCell scroll position.
$ws2.Activate()
$ws2.Range("A2").Select()
